# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Valor Mora" detail table (rows 16-25, 10 rows covering 5 employees)
# is replaced with an updated table of 5 rows covering 4 employees. Deleting
# rows 20-24 removes the surplus rows while letting the former last data row
# (25, which carries the table's special "closing" border style) slide up
# into row 20 - exactly the shape of the new table. The two footer rows
# (signature block) follow along from 30/31 down to 25/26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop 5 of the 10 detail rows so the table shrinks from 10 to 5 rows; the
# previously-last row (with the distinct bottom-border style) becomes the
# new row 20, and the footer (signature) rows shift from 30/31 to 25/26.
$ws.Range("20:24").EntireRow.Delete()

# Header summary figures
$ws.Range("E11").Value = 267760   # VALOR MORA total
$ws.Range("C13").Value = 4        # Cant. Trabajadores
$ws.Range("F13").Value = 4        # Cant. Periodos

# Row 16 - ISAAC BUSTILLO
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1001969"
$ws.Range("D16").Value = "ISAAC BUSTILLO"
$ws.Range("E16").Value = "2206"
$ws.Range("F16").Value = 40000
$ws.Range("G16").Value = 1000000

# Row 17 - SERGIO LUIS VERGARA CASTELLON, periodo 2506
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047419329"
$ws.Range("D17").Value = "SERGIO LUIS VERGARA CASTELLON"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Row 18 - SERGIO LUIS VERGARA CASTELLON, periodo 2507
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047419329"
$ws.Range("D18").Value = "SERGIO LUIS VERGARA CASTELLON"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Row 19 - JHON GLEICER CARDENAS ORTEGA, periodo 2507
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1001898521"
$ws.Range("D19").Value = "JHON GLEICER CARDENAS ORTEGA"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# Row 20 - DEIVIS ARLEY PAUTT SUAREZ, periodo 2508 (new employee)
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "73576525"
$ws.Range("D20").Value = "DEIVIS ARLEY PAUTT SUAREZ"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500
